$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.4
$ws.Range("K9").Value = 1.95
$ws.Range("L9").Value = 6.5
$ws.Range("Z9").Value = 12
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 26
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 21
$ws.Range("AN9").Value = 3.4
$ws.Range("AO9").Value = 9.5
$ws.Range("AQ9").Value = 34
$ws.Range("AU9").Value = 11
$ws.Range("AV9").Value = 101
$ws.Range("AW9").Value = 7
$ws.Range("BA9").Value = 251

# Row 19
$ws.Range("G19").Value = 1.55
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 5.75
$ws.Range("J19").Value = 2.1
$ws.Range("K19").Value = 2.38
$ws.Range("L19").Value = 5.5
$ws.Range("U19").Value = 1.8
$ws.Range("V19").Value = 1.91
$ws.Range("W19").Value = 7.5
$ws.Range("X19").Value = 7.5
$ws.Range("Z19").Value = 11
$ws.Range("AA19").Value = 12
$ws.Range("AC19").Value = 12
$ws.Range("AH19").Value = 17
$ws.Range("AI19").Value = 29
$ws.Range("AN19").Value = 3.6
$ws.Range("AP19").Value = 17
$ws.Range("AQ19").Value = 23
$ws.Range("AX19").Value = 29
$ws.Range("AY19").Value = 34
$ws.Range("AZ19").Value = 101

# Row 22
$ws.Range("G22").Value = 1.7
$ws.Range("H22").Value = 3.5
$ws.Range("I22").Value = 5.25
$ws.Range("J22").Value = 2.38
$ws.Range("L22").Value = 5.5
$ws.Range("N22").Value = 8.5
$ws.Range("AC22").Value = 8.5
$ws.Range("AE22").Value = 17
$ws.Range("AF22").Value = 51
$ws.Range("AG22").Value = 251
$ws.Range("AH22").Value = 12
$ws.Range("AI22").Value = 26
$ws.Range("AJ22").Value = 17
$ws.Range("AM22").Value = 51
$ws.Range("AO22").Value = 9
$ws.Range("AQ22").Value = 29
$ws.Range("AU22").Value = 9
$ws.Range("AX22").Value = 29

# Row 28
$ws.Range("K28").Value = 2.38
$ws.Range("W28").Value = 17
$ws.Range("AH28").Value = 6.5
$ws.Range("AK28").Value = 9
$ws.Range("AM28").Value = 29
$ws.Range("AN28").Value = 8.5
$ws.Range("AY28").Value = 19
$ws.Range("AZ28").Value = 19

# Row 51
$ws.Range("G51").Value = 1.65
$ws.Range("H51").Value = 3.7
$ws.Range("I51").Value = 5.5
$ws.Range("J51").Value = 2.3
$ws.Range("K51").Value = 2.1
$ws.Range("L51").Value = 5.5
$ws.Range("M51").Value = 1.07
$ws.Range("N51").Value = 8.5
$ws.Range("Q51").Value = 2.1
$ws.Range("R51").Value = 1.7
$ws.Range("Z51").Value = 12
$ws.Range("AB51").Value = 34
$ws.Range("AC51").Value = 8.5
$ws.Range("AD51").Value = 7
$ws.Range("AE51").Value = 19
$ws.Range("AG51").Value = 451
$ws.Range("AK51").Value = 51
$ws.Range("AN51").Value = 3.5
$ws.Range("AQ51").Value = 29
$ws.Range("AX51").Value = 29

# Row 78
$ws.Range("M78").Value = 1.11
$ws.Range("N78").Value = 6.5
$ws.Range("O78").Value = 1.5
$ws.Range("P78").Value = 2.5
$ws.Range("AG78").Value = 1000

# Row 80
$ws.Range("G80").Value = 1.5
$ws.Range("H80").Value = 3.7
$ws.Range("I80").Value = 6.5
$ws.Range("J80").Value = 2.02
$ws.Range("L80").Value = 6.2
$ws.Range("M80").Value = 1.06
$ws.Range("N80").Value = 9.369999999999999
$ws.Range("O80").Value = 1.33
$ws.Range("P80").Value = 2.8
$ws.Range("Q80").Value = 1.98
$ws.Range("R80").Value = 1.65
$ws.Range("S80").Value = 1.4
$ws.Range("T80").Value = 2.55
$ws.Range("U80").Value = 2.07
$ws.Range("V80").Value = 1.6
$ws.Range("W80").Value = 5.5
$ws.Range("X80").Value = 6.1
$ws.Range("Y80").Value = 8.5
$ws.Range("Z80").Value = 10
$ws.Range("AB80").Value = 35
$ws.Range("AC80").Value = 8.75
$ws.Range("AD80").Value = 7.5
$ws.Range("AE80").Value = 21
$ws.Range("AF80").Value = 120
$ws.Range("AH80").Value = 14.5
$ws.Range("AI80").Value = 40
$ws.Range("AJ80").Value = 21
$ws.Range("AK80").Value = 150
$ws.Range("AL80").Value = 80
$ws.Range("AM80").Value = 80
$ws.Range("AN80").Value = 3.15
$ws.Range("AO80").Value = 7
$ws.Range("AP80").Value = 18.5
$ws.Range("AQ80").Value = 22
$ws.Range("AS80").Value = 300
$ws.Range("AT80").Value = 2.5
$ws.Range("AU80").Value = 8.25
$ws.Range("AV80").Value = 90
$ws.Range("AW80").Value = 7.6
$ws.Range("AX80").Value = 40
$ws.Range("AY80").Value = 45
$ws.Range("AZ80").Value = 300
$ws.Range("BA80").Value = 300
